# Apply updated cryptocurrency price / volume figures to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "567.60") must be
# forced to remain text, otherwise Excel auto-converts them to a numeric
# value and the original text formatting (trailing zeros, etc.) is lost.
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $style

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.60"
$ws.Range("D5").Style = $style

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.34"
$ws.Range("D6").Style = $style

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.83"
$ws.Range("D10").Style = $style

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = $style

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("D19").Style = $style

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.63"
$ws.Range("D20").Style = $style

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.39"
$ws.Range("D21").Style = $style

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.38"
$ws.Range("D22").Style = $style

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.81"
$ws.Range("D24").Style = $style

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = $style

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("D28").Style = $style

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").Style = $style

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.71"
$ws.Range("D33").Style = $style

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.12"
$ws.Range("D34").Style = $style

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.09"
$ws.Range("D35").Style = $style

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.885"
$ws.Range("D36").Style = $style

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = $style

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.57"
$ws.Range("D39").Style = $style

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.55"
$ws.Range("D41").Style = $style

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = $style

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("D45").Style = $style

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0543"
$ws.Range("D46").Style = $style

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.35"
$ws.Range("D47").Style = $style

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.45"
$ws.Range("D48").Style = $style

# Remaining cells (already non-numeric text, or percentage strings) can be
# assigned directly.
$ws.Range("D2").Value = "60.480.65"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.623.49"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +6.13%  "
$ws.Range("E6").Value = "  +2.77%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("D9").Value = "2.643.11"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  +5.70%  "
$ws.Range("E12").Value = "  +7.31%  "
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "3.091.85"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "60.171.25"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("E16").Value = "  +6.76%  "
$ws.Range("E17").Value = "  +5.19%  "
$ws.Range("D18").Value = "2.639.78"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +5.24%  "
$ws.Range("D29").Value = "0.0₃0800"
$ws.Range("E29").Value = "  +10.57%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("E32").Value = "  +5.06%  "
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("E36").Value = "  +7.23%  "
$ws.Range("E37").Value = "  +5.50%  "
$ws.Range("E38").Value = "  +8.78%  "
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("E40").Value = "  +7.27%  "
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  +4.27%  "
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +5.99%  "
$ws.Range("E48").Value = "  +15.04%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("E51").Value = "  +7.18%  "
